# Daily attendance processing - 2025-12-03 06:37:18
# Re-order the "Recorded By" email lists for several sessions, update a
# few derived statistic counters, and flip session D17 (PARASITOLOGY /
# Session 5) from "Pending" over to "Not Recorded" (matching the same
# look as the other not-yet-happened sessions, e.g. row 11/13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recorded By reordering (same set of recipients, new order) ---
$ws.Range("G2").Value  = "Amira.Sobhy@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G3").Value  = "System, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G4").Value  = "hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G5").Value  = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G6").Value  = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G7").Value  = "AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G12").Value = "yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G30").Value = "yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"

# --- Updated roll-up counters (Missing Sessions / Pending Sessions) ---
$ws.Range("L7").Value  = 3
$ws.Range("L8").Value  = 7
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 7

# --- Row 17 (PARASITOLOGY / Session 5) flips from "Pending" to
#     "Not Recorded": copy the formatting already used for "Not Recorded"
#     rows (e.g. row 11) onto row 17, then update the status text. ---
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A17:I17").PasteSpecial(-4122) | Out-Null
$ws.Range("I17").Value = "Not Recorded"
